$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's last two data rows (146 "Provincia de Limarí" and 147
# "Provincia de Linares") are being replaced with a new week's price report;
# their previous values get pushed down and kept as history, and two more
# historical rows (old 148/149) shift down to make room.

# Step 1: insert two blank rows at 148 so the existing rows 148-149 move to
# rows 150-151 (plain OOXML row-shift, formatting carried along).
$ws.Rows("148:149").Insert()

# Step 2: the current (not yet overwritten) rows 146 and 147 become the new
# rows 148 and 149 - copy them down before they get overwritten in step 3.
$ws.Range("A146:T146").Copy()
$ws.Range("A148:T148").PasteSpecial()
$ws.Range("A147:T147").Copy()
$ws.Range("A149:T149").PasteSpecial()

# Step 3: write the new weekly price data into row 146.
$ws.Range("D146").Value = 44448
$ws.Range("M146").Value = 240
$ws.Range("N146").Value = 15000
$ws.Range("O146").Value = 15000
$ws.Range("P146").Value = 15000
$ws.Range("Q146").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R146").Value = "Provincia del Elquí"
$ws.Range("S146").Value = 10000
$ws.Range("T146").Value = 1.5

# Step 4: write the new weekly price data into row 147.
$ws.Range("D147").Value = 44448
$ws.Range("L147").Value = "Segunda"
$ws.Range("M147").Value = 1250
$ws.Range("N147").Value = 11000
$ws.Range("O147").Value = 11000
$ws.Range("P147").Value = 11000
$ws.Range("S147").Value = 5500
